$wb = $excel.ActiveWorkbook

# --- Overview sheet: localization status strings now reflect a completed handback ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Columns.Item(5).ColumnWidth = 29.17
$ovw.Columns.Item(6).ColumnWidth = 29.17

$mdName = "96555b91-741f-48b4-9887-5c2f343ac0d9.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/96555b91-741f-48b4-9887-5c2f343ac0d9.md"

# --- zh-cn handback report ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("I2").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdName)
$zh.Range("J2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.c39b27a6ee30c08d8156d4f335606b59bf26ae0d.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-25 09:03:33"
$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# --- de-de handback report ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("I2").Value = $mdName
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdName)
$de.Range("J2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.c39b27a6ee30c08d8156d4f335606b59bf26ae0d.de-de.xlf"
$de.Range("K2").Value = "2016-08-25 09:03:39"
$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17
